# Regenerate save_data to use K (strikeouts) instead of Strike# (strike pitch count)
# in column G of Sheet1, rows 2-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 4
    3  = 9
    4  = 8
    5  = 4
    6  = 7
    7  = 1
    8  = 10
    9  = 9
    10 = 12
    11 = 9
    12 = 8
    13 = 6
    14 = 8
    15 = 10
    16 = 10
    17 = 9
    18 = 10
    19 = 5
    20 = 14
    21 = 7
    22 = 5
    23 = 3
    24 = 8
    25 = 5
    26 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
